$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Long text values used below (here-strings preserve line breaks and
# any trailing spaces / accented characters verbatim).
# ------------------------------------------------------------------
$newProf = @"
519033 - Carlos Yujiro Shigue
"@

$objectives = @"
Apresentar noções de mecânica dos fluidos e reologia, mediante estudo dos meios fluidos quando estáticos ou em movimento. Capacitar o aluno a modelar e resolver problemas de interesse em mecânica dos fluidos e reologia, com escolha adequada de hipóteses e aplicação de ferramentas correspondentes de solução.
"@

$shortSyllabus = @"
Fundamentos de mecânica dos fluidos. Revisão de estática dos fluidos. Formulação integral e diferencial das equações de transporte de massa, energia e quantidade de movimento. Análise dimensional e semelhança. Escoamento incompressível de fluidos ideais e viscosos, regime laminar e turbulento. Equação de Navier-Stokes. Teoria da camada limite. Escoamento de fluidos não newtonianos. Formulação tensorial: tensão e deformação. Viscosidade e reometria. Viscoelasticidade. Aplicações.
"@

$syllabus = @"
Introdução: conceito de fluido; propriedades e conceito de contínuo; modelagem de processos de transferência; métodos de análise; dimensões e unidades.
Revisão de estática de fluidos: equação básica da hidrostática, variação de pressão em um fluido estático; princípios de Stevin, de Pascal e de Arquimedes.
Formulação integral das equações de transporte: teorema de transporte de Reynolds; aplicação para os princípios de conservação de massa, quantidade de movimento e energia; equação de Bernoulli.
Formulação diferencial das equações de transporte: descrição do escoamento; forma diferencial: dos princípios de conservação de massa, quantidade de movimento e energia; formulação adimensional, análise dimensional e semelhança. Grupos adimensionais: número de Reynolds e número de Grashoff.
Escoamento incompressível interno: equações de Euler; lei de Newton para a viscosidade, tensões de cisalhamento; equação de Navier-Stokes; regimes de escoamento: escoamento laminar e turbulento. Cálculo de perda de carga (distribuída e localizada), coeficiente de atrito. 
Escoamento incompressível externo: introdução à camada limite; escoamento ao redor de corpos, força da arraste.
Introdução a reologia. Definição e formulação tensorial de tensão e deformação. Tipos de deformação e escoamento de materiais. Equações fundamentais da reologia. Escoamento de fluidos newtonianos e não newtonianos. Viscosimetria e reometria. Reologia de sistemas dispersos. Colóides e emulsões. Soluções diluídas. Viscosimetria capilar. Aplicações.
"@

$method = @"
A avaliação será feita por meio de duas provas escritas P1 e P2 e por listas de exercícios e relatórios.
"@

$criteria = @"
A Nota final (NF) será calculada pela média ponderada das provas escritas e pela média dos trabalhos TR da seguinte maneira: NF = (P1 + 2*P2 + TR)/4
"@

$recoveryNorm = @"
A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2
"@

$bibliography = @"
BIRD,R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. LTC Editora, 2004.
FOX, R. W., McDONALD, A. T. Introdução à Mecânica dos Fluidos. LTC Editora, 2001.
SISSOM, L. E., PITTS, D. R. Fenômenos de Transporte. Ed. Guanabara, 1988.
SCHRAMM, G. Reologia e Reometria. Editora Artliber, 2006.
MANRICH, S.; PESSAN, L.A. Reologia: Conceitos Básicos, Editora UFSCar, 1987.
MALKIN, A. Rheology Fundamentals. ChemTec Publishing, 1994.
"@

# ------------------------------------------------------------------
# 1) Insert a new row at 13 for "Docentes responsaveis" (professor)
#    value, which previously lived (mis-placed) in the "Objetivos"
#    row. Column A of the new row must stay empty (no label there),
#    so clear it after the native row-insert copies formatting down,
#    then paste the correct column B / C formatting onto B13 / C13.
# ------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()

$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B13").Value = $newProf
$ws.Range("C13").Value = $newProf

# ------------------------------------------------------------------
# 2) Replace the misplaced content with the real course text in each
#    of the affected rows (all rows below 13 shifted down by one).
# ------------------------------------------------------------------
$ws.Range("B10").Value = $objectives        # Objetivos:
$ws.Range("C10").Value = $objectives

$ws.Range("B14").Value = $shortSyllabus      # Programa resumido:
$ws.Range("C14").Value = $shortSyllabus

$ws.Range("B16").Value = $syllabus           # Programa:
$ws.Range("C16").Value = $syllabus

$ws.Range("B19").Value = $method             # Metodo:
$ws.Range("C19").Value = $method

$ws.Range("B20").Value = $criteria           # Criterio:
$ws.Range("C20").Value = $criteria

$ws.Range("B21").Value = $recoveryNorm       # Norma de recuperacao:
$ws.Range("C21").Value = $recoveryNorm

$ws.Range("B22").Value = $bibliography       # Bibliografia:
$ws.Range("C22").Value = $bibliography

# ------------------------------------------------------------------
# 3) Column A no longer needs to span columns 1-2 (column B already
#    carries its own, more specific width/style definition that takes
#    precedence) - touching column B's width forces the engine to give
#    column A its own single-column definition instead of the old
#    1-2 combined range, without altering column A's stored width.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

Write-Output "edit complete"
